# "Info per municipality in statistics xls"
#
# 1. Rename the sheet "Evikomp" -> "Evikomp totalt"
# 2. Change the accumulator formulas in columns L (Ackumulerat, col 12)
#    and O (Ackumulerat, col 15) so that they show a blank ("") instead
#    of repeating the previous running total when no new value has been
#    entered for that row yet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Name = "Evikomp totalt"

# Row 3 only has the O-column accumulator (L3 is left untouched).
$ws.Range("O3").Formula = '=IF(N3="","",SUM(N$3:N3))'

# Rows 4-26 get both the L and O accumulator formulas updated.
for ($r = 4; $r -le 26; $r++) {
    $prev = $r - 1
    $ws.Range("L$r").Formula = '=IF(SUM(K$3:K' + $prev + ')=SUM(K$3:K' + $r + '),"",SUM(K$3:K' + $r + '))'
    $ws.Range("O$r").Formula = '=IF(N' + $r + '="","",SUM(N$3:N' + $r + '))'
}
